$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a new row (row 2) below the existing "Meta"/"Venda" header row,
# with the text value "4000.0" in both A2 and B2. The leading apostrophe
# forces Excel to store the numeric-looking text as a literal string
# rather than converting it to the number 4000.
$ws.Cells.Item(2, 1).Value = "'4000.0"
$ws.Cells.Item(2, 2).Value = "'4000.0"
